$wb = $excel.ActiveWorkbook

# ALC row 9
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 7182.8667
$ws.Range("I9").Value = 7182.8667
$ws.Range("K9").Value = 7182.8667
$ws.Range("M9").Value = -7013.8667

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2124.625
$ws.Range("I137").Value = 1544.3334
$ws.Range("K137").Value = 4633.0002
$ws.Range("M137").Value = -2083.0002

# ARM row 22
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 2624.5
$ws.Range("I22").Value = 2624.5
$ws.Range("K22").Value = 2624.5
$ws.Range("M22").Value = -2325.5

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2535.0435
$ws.Range("I32").Value = 2365.8635
$ws.Range("K32").Value = 2365.8635
$ws.Range("M32").Value = -2078.8635

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 25003030
$ws.Range("I74").Value = 30305788
$ws.Range("J74").Value = 4320.2856
$ws.Range("K74").Value = 30305788
$ws.Range("L74").Value = 4320.2856
$ws.Range("M74").Value = -30304914
$ws.Range("N74").Value = -6068.2856

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 25003030
$ws.Range("I77").Value = 30305788
$ws.Range("J77").Value = 4320.2856
$ws.Range("K77").Value = 151528940
$ws.Range("L77").Value = 21601.428
$ws.Range("M77").Value = -151524572
$ws.Range("N77").Value = -30337.428

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5280182.5
$ws.Range("I132").Value = 3032689.8
$ws.Range("K132").Value = 9098069.399999999
$ws.Range("M132").Value = -9095539.399999999

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 11980976
$ws.Range("I134").Value = 12879254
$ws.Range("K134").Value = 38637762
$ws.Range("M134").Value = -38635227

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1794.7213
$ws.Range("I31").Value = 1842.44
$ws.Range("J31").Value = 1577.8182
$ws.Range("K31").Value = 1842.44
$ws.Range("L31").Value = 1577.8182
$ws.Range("M31").Value = -1547.44
$ws.Range("N31").Value = -2167.8182

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1794.7213
$ws.Range("I34").Value = 1842.44
$ws.Range("J34").Value = 1577.8182
$ws.Range("K34").Value = 1842.44
$ws.Range("L34").Value = 1577.8182
$ws.Range("M34").Value = -1640.44
$ws.Range("N34").Value = -1981.8182

# CRP row 68
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 96000
$ws.Range("J68").Value = 99000
$ws.Range("L68").Value = 99000
$ws.Range("N68").Value = -100498

# CRP row 71
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 96000
$ws.Range("J71").Value = 99000
$ws.Range("L71").Value = 297000
$ws.Range("N71").Value = -304488

# CRP row 95
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 16174.667
$ws.Range("J95").Value = 16174.667
$ws.Range("L95").Value = 16174.667
$ws.Range("N95").Value = -21666.667

# CRP row 96
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H96").Value = 25541.334
$ws.Range("J96").Value = 25541.334
$ws.Range("L96").Value = 25541.334
$ws.Range("N96").Value = -31033.334

# CRP row 97
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 53176
$ws.Range("J97").Value = 53176
$ws.Range("L97").Value = 53176
$ws.Range("N97").Value = -55158

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 10003724
$ws.Range("I134").Value = 14709145
$ws.Range("J134").Value = 4704.625
$ws.Range("K134").Value = 44127435
$ws.Range("L134").Value = 14113.875
$ws.Range("M134").Value = -44124900
$ws.Range("N134").Value = -19183.875

# CUL row 3
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 10212.556
$ws.Range("I3").Value = 9614.125
$ws.Range("K3").Value = 28842.375
$ws.Range("M3").Value = -28730.375

# CUL row 18
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1325.6
$ws.Range("I18").Value = 407
$ws.Range("K18").Value = 1221
$ws.Range("M18").Value = -1052

# CUL row 94
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 17911.273
$ws.Range("I94").Value = 5500
$ws.Range("K94").Value = 16500
$ws.Range("M94").Value = -15824

# CUL row 99
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

# CUL row 103
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 359.4
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1853.9
$ws.Range("I131").Value = 1256.5
$ws.Range("K131").Value = 3769.5
$ws.Range("M131").Value = 1270.5

# CUL row 136
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 3117.7144
$ws.Range("I136").Value = 3117.7144
$ws.Range("K136").Value = 9353.143199999999
$ws.Range("M136").Value = -4253.143199999999

# CUL row 138
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 3298
$ws.Range("I138").Value = 3298
$ws.Range("K138").Value = 9894
$ws.Range("M138").Value = -4754

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2643.2222
$ws.Range("I140").Value = 2643.2222
$ws.Range("K140").Value = 7929.6666
$ws.Range("M140").Value = -2749.6666

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2056009.2
$ws.Range("I132").Value = 2458485
$ws.Range("J132").Value = 3382.6
$ws.Range("K132").Value = 7375455
$ws.Range("L132").Value = 10147.8
$ws.Range("M132").Value = -7372925
$ws.Range("N132").Value = -15207.8

# GSM row 136
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 40837.668
$ws.Range("J136").Value = 40837.668
$ws.Range("L136").Value = 122513.004
$ws.Range("N136").Value = -127613.004

# LTW row 117
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H117").Value = 49999
$ws.Range("J117").Value = 49999
$ws.Range("L117").Value = 49999
$ws.Range("N117").Value = -59177

# WVR row 54
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 34998.2
$ws.Range("J54").Value = 34998.2
$ws.Range("L54").Value = 34998.2
$ws.Range("N54").Value = -36038.2

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2027.8125
$ws.Range("I107").Value = 1280
$ws.Range("J107").Value = 2367.7273
$ws.Range("K107").Value = 3840
$ws.Range("L107").Value = 7103.1819
$ws.Range("M107").Value = -1920
$ws.Range("N107").Value = -10943.1819

# WVR row 116
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H116").Value = 124999.5
$ws.Range("J116").Value = 124999.5
$ws.Range("L116").Value = 124999.5
$ws.Range("N116").Value = -134177.5

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 997.53845
$ws.Range("I126").Value = 915.2727
$ws.Range("J126").Value = 1450
$ws.Range("K126").Value = 2745.8181
$ws.Range("L126").Value = 4350
$ws.Range("M126").Value = -275.8181
$ws.Range("N126").Value = -9290

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 20836098
$ws.Range("I136").Value = 22730118
$ws.Range("J136").Value = 1887.5
$ws.Range("K136").Value = 68190354
$ws.Range("L136").Value = 5662.5
$ws.Range("M136").Value = -68187804
$ws.Range("N136").Value = -10762.5
